$wb2 = $excel.ActiveWorkbook
$ws = $wb2.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.984.05"
$ws.Range("E2").Value = "  +1.88%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.763.95"
$ws.Range("E3").Value = "  +0.48%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  -0.30%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "322.22"
$ws.Range("E5").Value = "  -0.24%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9982"
$ws.Range("E6").Value = "  -0.20%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4246"
$ws.Range("E7").Value = "  -3.97%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3598"
$ws.Range("E8").Value = "  -2.82%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "44.26"
$ws.Range("E9").Value = "  -1.30%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07468"
$ws.Range("E10").Value = "  -3.13%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.106"
$ws.Range("E11").Value = "  -0.31%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.9987"
$ws.Range("E12").Value = "  -0.18%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "21.47"
$ws.Range("E13").Value = "  -0.22%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.102"
$ws.Range("E14").Value = "  -0.69%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.355"
$ws.Range("E15").Value = "  -0.63%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.795.37"
$ws.Range("E16").Value = "  +1.74%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "91.71"
$ws.Range("E17").Value = "  +1.83%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001061"
$ws.Range("E18").Value = "  -0.86%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06387"
$ws.Range("E19").Value = "  +2.47%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.9974"
$ws.Range("E20").Value = "  -0.30%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.15"
$ws.Range("E21").Value = "  -0.98%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.984"
$ws.Range("E22").Value = "  -2.78%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "27.996.02"
$ws.Range("E23").Value = "  +1.75%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.31"
$ws.Range("E24").Value = "  -1.57%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.160"
$ws.Range("E25").Value = "  -6.20%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "158.77"
$ws.Range("E26").Value = "  +4.05%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.21"
$ws.Range("E27").Value = "  -1.17%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.991.88"
$ws.Range("E28").Value = "  +1.78%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.147"
$ws.Range("E29").Value = "  -5.93%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "126.19"
$ws.Range("E30").Value = "  -0.63%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.180"
$ws.Range("E31").Value = "  +1.01%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.686"
$ws.Range("E32").Value = "  -0.14%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09076"
$ws.Range("E33").Value = "  -1.13%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.507"
$ws.Range("E34").Value = "  -3.22%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "12.60"
$ws.Range("E35").Value = "  +0.28%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02334"
$ws.Range("E36").Value = "  +1.45%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.060"
$ws.Range("E37").Value = "  +0.53%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2115"
$ws.Range("E38").Value = "  -1.60%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06080"
$ws.Range("E39").Value = "  -0.26%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.6414"
$ws.Range("E40").Value = "  -0.17%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.187"
$ws.Range("E41").Value = "  +0.79%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9974"
$ws.Range("E42").Value = "  -0.28%  "

$ws.Range("B43").Value = "WEMIXTOKEN"
$ws.Range("C43").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.395"
$ws.Range("E43").Value = "  +1.06%  "

$ws.Range("B44").Value = "FraxShare"
$ws.Range("C44").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "7.828"
$ws.Range("E44").Value = "  -0.94%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.60"
$ws.Range("E45").Value = "  +0.35%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5963"
$ws.Range("E46").Value = "  +0.39%  "

$ws.Range("E47").Value = "  -0.31%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.010"
$ws.Range("E48").Value = "  +1.96%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "123.42"
$ws.Range("E49").Value = "  -1.93%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.174"
$ws.Range("E50").Value = "  +3.84%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06885"
$ws.Range("E51").Value = "  +0.28%  "
